# updated fuel mix and nordnorgebanen cost
#
# On the "init_fuel_mix" sheet, the "Rail" / "Catenary" share (D9) and the
# "Rail" / "Methanol" share (D10) are updated, and the sheet's active
# selection moves from the stale J15 to D11 (matching what Excel leaves
# behind after editing D9/D10 and pressing Enter twice).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("init_fuel_mix")
$ws.Activate()

$ws.Range("D9").Value = 80
$ws.Range("D10").Value = 20

$ws.Range("D11").Select()
